$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SSC360 URL (B4) to the new CRM address and (re)create its hyperlink.
$newUrl = "https://sscpp-ppspc-360-ent.crm3.dynamics.com/"
$ws.Hyperlinks.Add($ws.Range("B4"), $newUrl, "", $newUrl, $newUrl)

# Match the existing "Hyperlink" cell style used elsewhere (e.g. D4).
$ws.Range("B4").Style = $ws.Range("D4").Style

# Widen column B to fit the longer URL text.
$ws.Columns("B").ColumnWidth = 89.15

# Update the saved selection/active cell.
$ws.Range("B15").Select()
